$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 4140
$ws.Range("I6").Value = 233.33333
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 699.99999
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = -587.99999
$ws.Range("N6").Value = -30224
# Row 47
$ws.Range("H47").Value = 2935
$ws.Range("J47").Value = 1990.7407
$ws.Range("L47").Value = 1990.7407
$ws.Range("N47").Value = -3934.7407
# Row 86
$ws.Range("H86").Value = 2648.5
$ws.Range("I86").Value = 2426.4285
$ws.Range("J86").Value = 3166.6667
$ws.Range("K86").Value = 2426.4285
$ws.Range("L86").Value = 3166.6667
$ws.Range("M86").Value = -1303.4285
$ws.Range("N86").Value = -5412.6667
# Row 89
$ws.Range("H89").Value = 2648.5
$ws.Range("I89").Value = 2426.4285
$ws.Range("J89").Value = 3166.6667
$ws.Range("K89").Value = 12132.1425
$ws.Range("L89").Value = 15833.3335
$ws.Range("M89").Value = -6516.1425
$ws.Range("N89").Value = -27065.3335
# Row 125
$ws.Range("H125").Value = 53711.875
$ws.Range("I125").Value = 71032.664
$ws.Range("J125").Value = 1749.5
$ws.Range("K125").Value = 639293.976
$ws.Range("L125").Value = 15745.5
$ws.Range("M125").Value = -636833.976
$ws.Range("N125").Value = -20665.5
# Row 138
$ws.Range("H138").Value = 1328.9
$ws.Range("I138").Value = 1328.9
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 3986.7
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1153.3
$ws.Range("N138").ClearContents()

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3818696
$ws.Range("I32").Value = 732359.4
$ws.Range("K32").Value = 732359.4
$ws.Range("M32").Value = -732072.4
# Row 97
$ws.Range("H97").Value = 29169.184
$ws.Range("I97").Value = 7780
$ws.Range("J97").Value = 98089.89
$ws.Range("K97").Value = 7780
$ws.Range("L97").Value = 98089.89
$ws.Range("M97").Value = -7284
$ws.Range("N97").Value = -99081.89

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3259.1365
$ws.Range("I20").Value = 2259.1765
$ws.Range("J20").Value = 6659
$ws.Range("K20").Value = 2259.1765
$ws.Range("L20").Value = 6659
$ws.Range("M20").Value = -2012.1765
$ws.Range("N20").Value = -7153
# Row 105
$ws.Range("H105").Value = 16382.656
$ws.Range("J105").Value = 71881.336
$ws.Range("L105").Value = 71881.336
$ws.Range("N105").Value = -75375.336
# Row 107
$ws.Range("H107").Value = 2406
$ws.Range("I107").Value = 2400.4827
$ws.Range("J107").Value = 2432.6667
$ws.Range("K107").Value = 2400.4827
$ws.Range("L107").Value = 2432.6667
$ws.Range("M107").Value = -480.4827
$ws.Range("N107").Value = -6272.6667

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3068.5386
$ws.Range("I16").Value = 2534.818
$ws.Range("K16").Value = 2534.818
$ws.Range("M16").Value = -2247.818
# Row 31
$ws.Range("H31").Value = 4278.385
$ws.Range("I31").Value = 2654.1428
$ws.Range("J31").Value = 6173.3335
$ws.Range("K31").Value = 2654.1428
$ws.Range("L31").Value = 6173.3335
$ws.Range("M31").Value = -2359.1428
$ws.Range("N31").Value = -6763.3335
# Row 34
$ws.Range("H34").Value = 4278.385
$ws.Range("I34").Value = 2654.1428
$ws.Range("J34").Value = 6173.3335
$ws.Range("K34").Value = 2654.1428
$ws.Range("L34").Value = 6173.3335
$ws.Range("M34").Value = -2452.1428
$ws.Range("N34").Value = -6577.3335
# Row 99
$ws.Range("H99").Value = 21158.945
$ws.Range("I99").Value = 25542.908
$ws.Range("K99").Value = 25542.908
$ws.Range("M99").Value = -24044.908
# Row 113
$ws.Range("H113").Value = 3068.5386
$ws.Range("I113").Value = 2534.818
$ws.Range("K113").Value = 2534.818
$ws.Range("M113").Value = -364.8180000000002
# Row 126
$ws.Range("H126").Value = 21158.945
$ws.Range("I126").Value = 25542.908
$ws.Range("K126").Value = 76628.724
$ws.Range("M126").Value = -74158.724
# Row 134
$ws.Range("H134").Value = 1208.5
$ws.Range("J134").Value = 1581
$ws.Range("L134").Value = 4743
$ws.Range("N134").Value = -9813

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 7361677.5
$ws.Range("I4").Value = 7084749.5
$ws.Range("K4").Value = 21254248.5
$ws.Range("M4").Value = -21254136.5
# Row 5
$ws.Range("H5").Value = 890
$ws.Range("I5").Value = 717.55554
$ws.Range("J5").Value = 1096.9333
$ws.Range("K5").Value = 2152.66662
$ws.Range("L5").Value = 3290.7999
$ws.Range("M5").Value = -2040.66662
$ws.Range("N5").Value = -3514.7999
# Row 14
$ws.Range("H14").Value = 260.15384
$ws.Range("I14").Value = 260.15384
$ws.Range("K14").Value = 780.4615200000001
$ws.Range("M14").Value = -607.4615200000001
# Row 35
$ws.Range("H35").Value = 427.5
$ws.Range("J35").Value = 925
$ws.Range("L35").Value = 2775
$ws.Range("N35").Value = -3351
# Row 68
$ws.Range("H68").Value = 1250.5
$ws.Range("I68").Value = 1298.5714
$ws.Range("K68").Value = 3895.7142
$ws.Range("M68").Value = -3084.7142
# Row 71
$ws.Range("H71").Value = 1250.5
$ws.Range("I71").Value = 1298.5714
$ws.Range("K71").Value = 11687.1426
$ws.Range("M71").Value = -7631.142600000001
# Row 129
$ws.Range("H129").Value = 119701.53
$ws.Range("J129").Value = 5134.6665
$ws.Range("L129").Value = 15403.9995
$ws.Range("N129").Value = -25403.9995
# Row 135
$ws.Range("H135").Value = 890
$ws.Range("I135").Value = 717.55554
$ws.Range("J135").Value = 1096.9333
$ws.Range("K135").Value = 6457.99986
$ws.Range("L135").Value = 9872.3997
$ws.Range("M135").Value = -3922.99986
$ws.Range("N135").Value = -14942.3997

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 735.6842
$ws.Range("I2").Value = 306.6154
$ws.Range("K2").Value = 306.6154
$ws.Range("M2").Value = -193.6154
# Row 70
$ws.Range("H70").Value = 54117030
$ws.Range("I70").Value = 5442.3076
$ws.Range("J70").Value = 288600580
$ws.Range("K70").Value = 5442.3076
$ws.Range("L70").Value = 288600580
$ws.Range("M70").Value = -5172.3076
$ws.Range("N70").Value = -288601120
# Row 73
$ws.Range("H73").Value = 54117030
$ws.Range("I73").Value = 5442.3076
$ws.Range("J73").Value = 288600580
$ws.Range("K73").Value = 5442.3076
$ws.Range("L73").Value = 288600580
$ws.Range("M73").Value = -4506.3076
$ws.Range("N73").Value = -288602452
# Row 80
$ws.Range("H80").Value = 9712.333000000001
$ws.Range("I80").Value = 14574.5
$ws.Range("J80").Value = 4155.5713
$ws.Range("K80").Value = 14574.5
$ws.Range("L80").Value = 4155.5713
$ws.Range("M80").Value = -13576.5
$ws.Range("N80").Value = -6151.5713
# Row 83
$ws.Range("H83").Value = 9712.333000000001
$ws.Range("I83").Value = 14574.5
$ws.Range("J83").Value = 4155.5713
$ws.Range("K83").Value = 72872.5
$ws.Range("L83").Value = 20777.8565
$ws.Range("M83").Value = -67880.5
$ws.Range("N83").Value = -30761.8565
# Row 97
$ws.Range("H97").Value = 69253.78999999999
$ws.Range("I97").Value = 41223.688
$ws.Range("K97").Value = 41223.688
$ws.Range("M97").Value = -40727.688
# Row 126
$ws.Range("H126").Value = 41408.438
$ws.Range("I126").Value = 2001.875
$ws.Range("K126").Value = 6005.625
$ws.Range("M126").Value = -3535.625
# Row 132
$ws.Range("H132").Value = 4704.7915
$ws.Range("I132").Value = 5005.875
$ws.Range("J132").Value = 4102.625
$ws.Range("K132").Value = 15017.625
$ws.Range("L132").Value = 12307.875
$ws.Range("M132").Value = -12487.625
$ws.Range("N132").Value = -17367.875

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 18421.25
$ws.Range("J6").Value = 21561.666
$ws.Range("L6").Value = 21561.666
$ws.Range("N6").Value = -21785.666
# Row 7
$ws.Range("H7").Value = 9341.208000000001
$ws.Range("I7").Value = 7532.4443
$ws.Range("J7").Value = 14767.5
$ws.Range("K7").Value = 7532.4443
$ws.Range("L7").Value = 14767.5
$ws.Range("M7").Value = -7420.4443
$ws.Range("N7").Value = -14991.5
# Row 40
$ws.Range("H40").Value = 6541
$ws.Range("I40").Value = 6395.6
$ws.Range("J40").Value = 7995
$ws.Range("K40").Value = 6395.6
$ws.Range("L40").Value = 7995
$ws.Range("M40").Value = -6259.6
$ws.Range("N40").Value = -8267
# Row 61
$ws.Range("H61").Value = 94016.94
$ws.Range("I61").Value = 88786.914
$ws.Range("J61").Value = 109707
$ws.Range("K61").Value = 88786.914
$ws.Range("L61").Value = 109707
$ws.Range("M61").Value = -88584.914
$ws.Range("N61").Value = -110111
# Row 113
$ws.Range("H113").Value = 94016.94
$ws.Range("I113").Value = 88786.914
$ws.Range("J113").Value = 109707
$ws.Range("K113").Value = 88786.914
$ws.Range("L113").Value = 109707
$ws.Range("M113").Value = -86616.914
$ws.Range("N113").Value = -114047
# Row 126
$ws.Range("H126").Value = 9341.208000000001
$ws.Range("I126").Value = 7532.4443
$ws.Range("J126").Value = 14767.5
$ws.Range("K126").Value = 22597.3329
$ws.Range("L126").Value = 44302.5
$ws.Range("M126").Value = -20127.3329
$ws.Range("N126").Value = -49242.5
# Row 132
$ws.Range("H132").Value = 5046
$ws.Range("I132").Value = 4878.9375
$ws.Range("K132").Value = 14636.8125
$ws.Range("M132").Value = -12106.8125

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 8314.817999999999
$ws.Range("I81").Value = 6292.7
$ws.Range("J81").Value = 9999.916999999999
$ws.Range("K81").Value = 12585.4
$ws.Range("L81").Value = 19999.834
$ws.Range("M81").Value = -11524.4
$ws.Range("N81").Value = -22121.834
# Row 84
$ws.Range("H84").Value = 8314.817999999999
$ws.Range("I84").Value = 6292.7
$ws.Range("J84").Value = 9999.916999999999
$ws.Range("K84").Value = 62927
$ws.Range("L84").Value = 99999.17
$ws.Range("M84").Value = -57623
$ws.Range("N84").Value = -110607.17
# Row 107
$ws.Range("H107").Value = 33372866
$ws.Range("I107").Value = 1861.5834
$ws.Range("K107").Value = 5584.7502
$ws.Range("M107").Value = -3664.7502
# Row 126
$ws.Range("H126").Value = 2348.3845
$ws.Range("I126").Value = 2012.3334
$ws.Range("K126").Value = 6037.0002
$ws.Range("M126").Value = -3567.0002

